$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70
$ws.Cells.Item(70, 1).Value = 3
$ws.Cells.Item(70, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 44628
$ws.Cells.Item(70, 5).Value = 5
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100104
$ws.Cells.Item(70, 8).Value = "Frutos de pepita"
$ws.Cells.Item(70, 9).Value = 100104003
$ws.Cells.Item(70, 10).Value = "Membrillo"
$ws.Cells.Item(70, 11).Value = "Champion"
$ws.Cells.Item(70, 12).Value = "Especial"
$ws.Cells.Item(70, 13).Value = 50
$ws.Cells.Item(70, 14).Value = 15000
$ws.Cells.Item(70, 15).Value = 15000
$ws.Cells.Item(70, 16).Value = 15000
$ws.Cells.Item(70, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(70, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(70, 19).Value = 833
$ws.Cells.Item(70, 20).Value = 18
$ws.Cells.Item(70, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 71
$ws.Cells.Item(71, 1).Value = 3
$ws.Cells.Item(71, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44628
$ws.Cells.Item(71, 5).Value = 5
$ws.Cells.Item(71, 6).Value = "Fruta"
$ws.Cells.Item(71, 7).Value = 100104
$ws.Cells.Item(71, 8).Value = "Frutos de pepita"
$ws.Cells.Item(71, 9).Value = 100104003
$ws.Cells.Item(71, 10).Value = "Membrillo"
$ws.Cells.Item(71, 11).Value = "Champion"
$ws.Cells.Item(71, 12).Value = "Primera"
$ws.Cells.Item(71, 13).Value = 58
$ws.Cells.Item(71, 14).Value = 13000
$ws.Cells.Item(71, 15).Value = 13000
$ws.Cells.Item(71, 16).Value = 13000
$ws.Cells.Item(71, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(71, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(71, 19).Value = 722
$ws.Cells.Item(71, 20).Value = 18
$ws.Cells.Item(71, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 72
$ws.Cells.Item(72, 1).Value = 3
$ws.Cells.Item(72, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(72, 3).Value = "Coquimbo"
$ws.Cells.Item(72, 4).Value = 44628
$ws.Cells.Item(72, 5).Value = 5
$ws.Cells.Item(72, 6).Value = "Fruta"
$ws.Cells.Item(72, 7).Value = 100104
$ws.Cells.Item(72, 8).Value = "Frutos de pepita"
$ws.Cells.Item(72, 9).Value = 100104003
$ws.Cells.Item(72, 10).Value = "Membrillo"
$ws.Cells.Item(72, 11).Value = "Champion"
$ws.Cells.Item(72, 12).Value = "Segunda"
$ws.Cells.Item(72, 13).Value = 56
$ws.Cells.Item(72, 14).Value = 11000
$ws.Cells.Item(72, 15).Value = 11000
$ws.Cells.Item(72, 16).Value = 11000
$ws.Cells.Item(72, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(72, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(72, 19).Value = 611
$ws.Cells.Item(72, 20).Value = 18
$ws.Cells.Item(72, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

